# Weekly data refresh: insert a new observation row at row 32
# (pushing all subsequent rows down by one) for the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Cebollín" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 32, shifting rows 32:142
# down to 33:143 (mirrors Excel's native Insert behaviour, which also
# copies the row-above formatting, e.g. the date style on column D).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly record.
$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(32, 3).Value = "Ñuble"
$ws.Cells.Item(32, 4).Value = 45070
$ws.Cells.Item(32, 5).Value = 16
$ws.Cells.Item(32, 6).Value = 100112037
$ws.Cells.Item(32, 7).Value = "Cebollín"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 80
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 7000
$ws.Cells.Item(32, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(32, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(32, 16).Value = 194
$ws.Cells.Item(32, 17).Value = 36
$ws.Cells.Item(32, 18).Value = "Hortaliza"
